$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 124, shifting existing rows 124:146 down to 125:147
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new weekly data point
$ws.Range("A124").Value = 11
$ws.Range("B124").Value = "Vega Monumental Concepción"
$ws.Range("C124").Value = "Bíobío"
$ws.Range("D124").Value = "2023-03-30"
$ws.Range("E124").Value = 8
$ws.Range("F124").Value = 100112001
$ws.Range("G124").Value = "Berenjena"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 220
$ws.Range("K124").Value = 7000
$ws.Range("L124").Value = 7500
$ws.Range("M124").Value = 7227
$ws.Range("N124").Value = "$/caja 60 unidades"
$ws.Range("O124").Value = "Región de Arica y Parinacota"
$ws.Range("P124").Value = 120
$ws.Range("Q124").Value = 60
$ws.Range("R124").Value = "Hortaliza"
